# "Updating static tables to 2023."
#
# The "unit_file" sheet (manual_corrections.xlsx -> unit_file) had a
# duplicate/obsolete "delete" correction row for plant_id 2132 (row 8).
# That row is removed entirely, which shifts every following row up by
# one (old row 9 -> new row 8, ... old row 25 -> new row 24).

$wb = $excel.ActiveWorkbook

$unit = $wb.Worksheets.Item("unit_file")

# Remove the obsolete correction row (plant_id 2132, column_to_update
# "delete"); this shifts all subsequent rows up by one, matching the new
# last row 24 (was 25).
$unit.Rows.Item(8).Delete() | Out-Null

# Mirror the author's final on-screen state: plant_file was the
# previously active sheet (selection left at A12) and unit_file becomes
# the active sheet again, selected at D12.
$plant = $wb.Worksheets.Item("plant_file")
$plant.Activate() | Out-Null
$plant.Range("A12").Select() | Out-Null

$unit.Activate() | Out-Null
$unit.Range("D12").Select() | Out-Null
